# Applies the commit "3.3 files and starter package" edit to
# BAU Dispatch Priority by Elec Source.xlsx
#
# Summary of the functional changes:
#  - About sheet: add "Colorado" label next to the title, bump the
#    date stamp in C1 from 4/21/2021 to 9/24/2021
#  - BDPbES sheet becomes the active/selected sheet (tab) instead of About
#  - BDPbES sheet: "onshore wind" (row 6) and "solar PV" (row 7) priority
#    values change from 2 to 1 across the whole year range (B:AK)

$wb = $excel.ActiveWorkbook

$wsAbout  = $wb.Worksheets.Item("About")
$wsBDPbES = $wb.Worksheets.Item("BDPbES")

# --- About sheet updates ---------------------------------------------------
$wsAbout.Range("B1").Value = "Colorado"
$wsAbout.Range("C1").Value = 44463

# --- BDPbES sheet updates ---------------------------------------------------
# Row 6 = onshore wind, Row 7 = solar PV -- set base (year 2015) value;
# the remaining year columns are driven by shared formulas ("=$B6"/"=$B7")
# so they recompute automatically.
$wsBDPbES.Range("B6").Value = 1
$wsBDPbES.Range("B7").Value = 1

$excel.Calculate()

# --- Make BDPbES the active/selected sheet ---------------------------------
$wsBDPbES.Activate()
$wsBDPbES.Range("B8").Select()

$wb.Save()
